$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F5: "X" -> "Sign", with a new blue fill style (mirrors the new
#     "Sign" bit flag added to the opcode encoding table) ---
$ws.Range("F5").Value = "Sign"
$ws.Range("F5").Interior.Color = 12611584   # RGB(0,112,192) = FF0070C0

# --- New "Sign" legend table in I26:L29, mirroring the existing
#     "Ld/str" legend table in E26:H29 ---
$ws.Range("I26").Value = "Sign"

$ws.Range("J27").Value = "Bit"
$ws.Range("G27").Copy()
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("K27").Value = 0

$ws.Range("G28").Copy()
$ws.Range("K28").PasteSpecial(-4122)
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = "Unsigned"

$ws.Range("G29").Copy()
$ws.Range("K29").PasteSpecial(-4122)
$ws.Range("K29").Value = 1
$ws.Range("L29").Value = "Signed"

# --- keep the view roughly where the author left it ---
$ws.Range("L29").Select()
$excel.ActiveWindow.ScrollRow = 13
